$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

# Fill in the "mediciones" ICQ values (column E) for rows 11-21
$ws.Range("E11").Value = 2.95
$ws.Range("F11").Value = 0.0094

$ws.Range("E12").Value = 14.21
$ws.Range("E13").Value = 19.82
$ws.Range("E14").Value = 20.04
$ws.Range("E15").Value = 14.39
$ws.Range("E16").Value = 20.03
$ws.Range("E17").Value = 15
$ws.Range("E18").Value = 14.85
$ws.Range("E19").Value = 20.23
$ws.Range("E20").Value = 1.27
$ws.Range("E21").Value = 1.26

# Update the selected cell to match the author's final cursor position
$ws.Range("E24").Select()
